$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.375.19"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.633.66"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.42"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.15"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.634.06"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("E10").Value = "  +11.73%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.59"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").Value = "  +5.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.109.22"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.219.98"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.647.09"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.39"
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "369.33"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.40"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.82"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.89"
$ws.Range("E25").Value = "  +7.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.89"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.757.34"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("E29").Value = "  +4.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "572.71"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.40"
$ws.Range("E38").Value = "  +2.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.15"
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.90"
$ws.Range("E40").Value = "  +5.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.367"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.34"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("E44").Value = "  +4.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0327"
$ws.Range("E45").Value = "  +11.01%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.37"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "155.08"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.70"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.92"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.70"
$ws.Range("E51").Value = "  +0.33%  "
